$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.067.87"
$ws.Range("E2").Value = "  +2.41%  "

$ws.Range("D3").Value = "1.654.15"
$ws.Range("E3").Value = "  +3.28%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.22%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +1.68%  "

$ws.Range("E9").Value = "  +1.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0874"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.98%  "

$ws.Range("D12").Value = "1.888.42"
$ws.Range("E12").Value = "  +3.32%  "

$ws.Range("D13").Value = "1.653.71"
$ws.Range("E13").Value = "  +3.24%  "

$ws.Range("E14").Value = "  +1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.47%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "27.065.04"
$ws.Range("E17").Value = "  +2.47%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("D20").Value = "0.0₃0730"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.02%  "

$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.87%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.114"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("E31").Value = "  +1.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.06%  "

$ws.Range("D33").Value = "1.514.01"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("E34").Value = "  +4.61%  "

$ws.Range("E35").Value = "  +8.65%  "

$ws.Range("E36").Value = "  -0.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.578"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.24%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.888"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.14%  "

$ws.Range("E40").Value = "  +2.68%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.11%  "

$ws.Range("D44").Value = "1.795.24"
$ws.Range("E44").Value = "  +3.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.775"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.917"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  +10.38%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.49%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0975"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
